$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.708.13"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.639.54"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'212.58"
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'23.21"
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "'0.0890"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "1.870.76"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "1.642.22"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E15").Value = "  -4.53%  "
$ws.Range("D16").Value = "'64.78"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "27.667.19"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "'230.61"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  +2.22%  "
$ws.Range("D20").Value = "0.0₃0722"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'4.32"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'10.25"
$ws.Range("E23").Value = "  +5.05%  "
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("D25").Value = "'150.82"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "1.458.99"
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("E34").Value = "  -2.37%  "
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("D37").Value = "'0.567"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "'0.881"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("D40").Value = "'0.901"
$ws.Range("E40").Value = "  +10.12%  "
$ws.Range("D41").Value = "'69.32"
$ws.Range("E41").Value = "  +6.21%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").Value = "  -1.73%  "
$ws.Range("D44").Value = "'5.61"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("D47").Value = "1.780.59"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("E48").Value = "  +2.92%  "
$ws.Range("D49").Value = "'86.98"
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").Value = "'0.0995"
$ws.Range("E51").Value = "  +0.38%  "
